$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.193.23'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '1.911.37'
$ws.Range('E3').Value = '  +2.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.85%  '
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5063'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3922'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09339'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.142'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.98'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.406'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.98'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('D14').Value = '1.916.08'
$ws.Range('E14').Value = '  +2.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.325'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.32%  '
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.55'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06627'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.05'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.15%  '
$ws.Range('E21').Value = '  -0.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.233'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.89%  '
$ws.Range('D23').Value = '28.259.31'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  +2.72%  '
$ws.Range('E25').Value = '  +1.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.595'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.05%  '
$ws.Range('D27').Value = '2.136.54'
$ws.Range('E27').Value = '  +2.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.20'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '158.24'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.27'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('E31').Value = '  +4.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1074'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.673'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.609'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.698'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06701'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02439'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2216'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.40%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.246'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.281'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6537'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.56'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.032'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.82%  '
$ws.Range('E44').Value = '  -0.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6129'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.40'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.32%  '
$ws.Range('B47').Value = 'WEMIXTOKEN'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.299'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.58%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.723'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.78%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.033'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '122.47'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('E51').Value = '  -0.49%  '
